$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet tbl1 (sheet1): header rename, many updated stats, drop D2, add rows 14-23
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("tbl1")

$ws1.Range("C1").Value = "mean_stdev"

$ws1.Range("A2").Value = "Change in school COVID-19 case rate"
$ws1.Range("B2").Value = "342 (-6.25, 7.81)"
$ws1.Range("C2").Value = "1.16 (1.87)"
$ws1.Range("D2").Value = $null

$ws1.Range("B3").Value = "342 (-13.67, 14747.43)"
$ws1.Range("C3").Value = "826.27 (2375.13)"
$ws1.Range("D3").Value = "0.004 (0.937)"

$ws1.Range("B4").Value = "338 (0, 98.7)"
$ws1.Range("C4").Value = "1.64 (9.21)"
$ws1.Range("D4").Value = "-0.046 (0.401)"

$ws1.Range("B5").Value = "338 (0, 56)"
$ws1.Range("C5").Value = "4.04 (7.73)"
$ws1.Range("D5").Value = "0.017 (0.754)"

$ws1.Range("B6").Value = "338 (0, 99.5)"
$ws1.Range("C6").Value = "13.31 (21.69)"
$ws1.Range("D6").Value = "0.038 (0.489)"

$ws1.Range("B7").Value = "338 (0, 100)"
$ws1.Range("C7").Value = "24.97 (26.59)"
$ws1.Range("D7").Value = "0.066 (0.229)"

$ws1.Range("B8").Value = "338 (0, 9.7)"
$ws1.Range("C8").Value = "0.24 (0.7)"
$ws1.Range("D8").Value = "-0.041 (0.453)"

$ws1.Range("B9").Value = "338 (0, 2.6)"
$ws1.Range("C9").Value = "0.02 (0.17)"
$ws1.Range("D9").Value = "-0.047 (0.386)"

$ws1.Range("B10").Value = "338 (0, 23.8)"
$ws1.Range("C10").Value = "3.86 (3.09)"
$ws1.Range("D10").Value = "-0.004 (0.936)"

$ws1.Range("B11").Value = "338 (0, 100)"
$ws1.Range("C11").Value = "51.92 (32)"
$ws1.Range("D11").Value = "-0.07 (0.202)"

$ws1.Range("A12").Value = "Percent free and reduced lunch"
$ws1.Range("B12").Value = "331 (0.4, 100)"
$ws1.Range("C12").Value = "51.31 (28.3)"
$ws1.Range("D12").Value = "0.027 (0.63)"

$ws1.Range("B13").Value = "342 (0, 99.94)"
$ws1.Range("C13").Value = "51.31 (27.83)"
$ws1.Range("D13").Value = "0.068 (0.21)"

$ws1.Range("A14").Value = "Midwest"
$ws1.Range("B14").Value = "55 (-2.87, 4)"
$ws1.Range("C14").Value = "0.43 (1.33)"

$ws1.Range("A15").Value = "Northeast"
$ws1.Range("B15").Value = "94 (-2.93, 7.67)"
$ws1.Range("C15").Value = "1.47 (2.02)"

$ws1.Range("A16").Value = "South"
$ws1.Range("B16").Value = "120 (-6.25, 7.5)"
$ws1.Range("C16").Value = "1.68 (2.06)"

$ws1.Range("A17").Value = "West"
$ws1.Range("B17").Value = "73 (-1.16, 7.81)"
$ws1.Range("C17").Value = "0.46 (1.24)"

$ws1.Range("A18").Value = "City"
$ws1.Range("B18").Value = "83 (-2.6, 6.74)"
$ws1.Range("C18").Value = "1.06 (1.78)"

$ws1.Range("A19").Value = "Rural"
$ws1.Range("B19").Value = "92 (-6.25, 7.5)"
$ws1.Range("C19").Value = "1.17 (2)"

$ws1.Range("A20").Value = "Suburb"
$ws1.Range("B20").Value = "119 (-2.93, 7.81)"
$ws1.Range("C20").Value = "1.3 (1.93)"

$ws1.Range("A21").Value = "Town"
$ws1.Range("B21").Value = "48 (-2.87, 5.17)"
$ws1.Range("C21").Value = "0.97 (1.65)"

$ws1.Range("A22").Value = "Region"
$ws1.Range("A23").Value = "Locale"

# ---------------------------------------------------------------------------
# Sheet tbl2 (sheet2): updated B/C/D/E/F stats for every strategy row
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("tbl2")

$ws2.Range("B2").Value = "342 (0, 1)"
$ws2.Range("D2").Value = 1.193
$ws2.Range("E2").Value = 0.91
$ws2.Range("F2").Value = "0.283 (0.276)"

$ws2.Range("B3").Value = "342 (0, 1)"
$ws2.Range("D3").Value = 1.217
$ws2.Range("E3").Value = 0.997
$ws2.Range("F3").Value = "0.22 (0.358)"

$ws2.Range("B4").Value = "342 (0, 1)"
$ws2.Range("D4").Value = 1.261
$ws2.Range("E4").Value = 0.878
$ws2.Range("F4").Value = "0.383 (0.08)"

$ws2.Range("B5").Value = "342 (0, 1)"
$ws2.Range("C5").Value = "0.15 (0.35)"
$ws2.Range("D5").Value = 1.194
$ws2.Range("E5").Value = 0.958
$ws2.Range("F5").Value = "0.236 (0.454)"

$ws2.Range("B6").Value = "342 (0, 1)"
$ws2.Range("C6").Value = "0.41 (0.49)"
$ws2.Range("D6").Value = 1.257
$ws2.Range("E6").Value = 1.018
$ws2.Range("F6").Value = "0.24 (0.244)"

$ws2.Range("B7").Value = "342 (0, 1)"
$ws2.Range("D7").Value = 1.216
$ws2.Range("E7").Value = 1.032
$ws2.Range("F7").Value = "0.184 (0.383)"

$ws2.Range("B8").Value = "342 (0, 1)"
$ws2.Range("C8").Value = "0.32 (0.47)"
$ws2.Range("D8").Value = 1.236
$ws2.Range("E8").Value = 0.997
$ws2.Range("F8").Value = "0.239 (0.27)"

$ws2.Range("B9").Value = "342 (0, 1)"
$ws2.Range("C9").Value = "0.25 (0.44)"
$ws2.Range("D9").Value = 1.151
$ws2.Range("E9").Value = 1.186
$ws2.Range("F9").Value = "-0.035 (0.875)"

$ws2.Range("B10").Value = "342 (0, 1)"
$ws2.Range("C10").Value = "0.04 (0.21)"
$ws2.Range("D10").Value = 1.176
$ws2.Range("E10").Value = 0.802
$ws2.Range("F10").Value = "0.375 (0.4)"

$ws2.Range("B11").Value = "342 (0, 1)"
$ws2.Range("C11").Value = "0.33 (0.47)"
$ws2.Range("D11").Value = 1.316
$ws2.Range("E11").Value = 0.84
$ws2.Range("F11").Value = "0.476 (0.019)"

# ---------------------------------------------------------------------------
# Sheet tbl3 (sheet3): rows re-ranked with new NCES-derived missingness counts
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("tbl3")

$ws3.Range("A2").Value = "Region"
$ws3.Range("B2").Value = 100
$ws3.Range("A3").Value = "State"
$ws3.Range("B3").Value = 100
$ws3.Range("A4").Value = "Percent two or more races"
$ws3.Range("B4").Value = 77
$ws3.Range("A5").Value = "Percent Asian"
$ws3.Range("B5").Value = 43
$ws3.Range("A6").Value = "Percent White"
$ws3.Range("B6").Value = 37
$ws3.Range("A7").Value = "Percent free and reduced lunch"
$ws3.Range("B7").Value = 35
$ws3.Range("A8").Value = "SVI Overall Rank"
$ws3.Range("B8").Value = 31
$ws3.Range("A9").Value = "Percent Black or African American"
$ws3.Range("B9").Value = 23
$ws3.Range("A10").Value = "School level"
$ws3.Range("B10").Value = 21
$ws3.Range("A11").Value = "School enrollment"
$ws3.Range("B11").Value = 12
$ws3.Range("A12").Value = "Percent Hispanic or Latino"
$ws3.Range("B12").Value = 8
$ws3.Range("A13").Value = "Change in county COVID-19 case rate"
$ws3.Range("B13").Value = 6
$ws3.Range("A14").Value = "Percent American Indian/Alaska Native"
$ws3.Range("B14").Value = 6
$ws3.Range("A15").Value = "Percent Native Hawaiian or other Pacific Islander"
$ws3.Range("B15").Value = 1
$ws3.Range("A16").Value = "Locale"
$ws3.Range("B16").Value = 0

# ---------------------------------------------------------------------------
# Sheet tbl4 (sheet4): rows re-ranked with new percentages
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("tbl4")

$ws4.Range("A2").Value = "HVAC systems"
$ws4.Range("B2").Value = 100
$ws4.Range("A3").Value = "Contact tracing"
$ws4.Range("B3").Value = 94
$ws4.Range("A4").Value = "Cleaning"
$ws4.Range("B4").Value = 90
$ws4.Range("A5").Value = "Screening and testing for students"
$ws4.Range("B5").Value = 89
$ws4.Range("A6").Value = "Universal masking requirements"
$ws4.Range("B6").Value = 66
$ws4.Range("A7").Value = "Physical distancing"
$ws4.Range("B7").Value = 12
$ws4.Range("A8").Value = "Vaccination offered"
$ws4.Range("B8").Value = 11
$ws4.Range("A9").Value = "Staying home when sick"
$ws4.Range("B9").Value = 2
$ws4.Range("A10").Value = "HEPA filters"
$ws4.Range("B10").Value = 1
$ws4.Range("A11").Value = "Quarantining"
$ws4.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# Sheet tbl5 (sheet5): updated coefficients
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("tbl5")

$ws5.Range("B2").Value = "-0.08 (-0.69, 0.53)"
$ws5.Range("B3").Value = "-0.13 (-0.54, 0.27)"
$ws5.Range("B4").Value = "-0.38 (-0.83, 0.03)"
$ws5.Range("B5").Value = "-0.35 (-0.84, 0.2)"
$ws5.Range("B6").Value = "-0.33 (-0.71, 0.03)"
$ws5.Range("B7").Value = "-0.17 (-0.54, 0.22)"
$ws5.Range("B8").Value = "-0.29 (-0.72, 0.08)"
$ws5.Range("B9").Value = "-0.09 (-0.53, 0.31)"
$ws5.Range("B10").Value = "0.03 (-0.98, 0.96)"
$ws5.Range("B11").Value = "-0.38 (-0.78, 0.02)"

# ---------------------------------------------------------------------------
# Sheet tbl6 (sheet6): updated coefficients + new rows 14-21
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("tbl6")

$ws6.Range("B2").Value = "1.04 (0.32, 1.78)"
$ws6.Range("B3").Value = "0.33 (-0.44, 1.1)"
$ws6.Range("B4").Value = "0.14 (-0.35, 0.68)"
$ws6.Range("B5").Value = "-0.28 (-0.84, 0.3)"
$ws6.Range("B6").Value = "-0.31 (-0.99, 0.36)"
$ws6.Range("B7").Value = "-0.2 (-0.85, 0.47)"
$ws6.Range("B8").Value = "0.13 (-0.42, 0.69)"
$ws6.Range("B9").Value = "-0.06 (-0.7, 0.56)"
$ws6.Range("B10").Value = "0.33 (-0.23, 0.82)"
$ws6.Range("B11").Value = "0.25 (-0.6, 1.29)"
$ws6.Range("B12").Value = "-0.4 (-1.01, 0.17)"
$ws6.Range("B13").Value = "0.24 (0.04, 0.46)"

$ws6.Range("A14").Value = "Percent Asian"
$ws6.Range("B14").Value = "0.01 (-0.21, 0.23)"

$ws6.Range("A15").Value = "Percent White"
$ws6.Range("B15").Value = "-0.15 (-0.5, 0.26)"

$ws6.Range("A16").Value = "Percent free and reduced lunch"
$ws6.Range("B16").Value = "-0.25 (-0.57, 0.1)"

$ws6.Range("A17").Value = "SVI Overall Rank"
$ws6.Range("B17").Value = "0.04 (-0.2, 0.28)"

$ws6.Range("A18").Value = "Percent Black or African American"
$ws6.Range("B18").Value = "-0.07 (-0.35, 0.18)"

$ws6.Range("A19").Value = "High school"
$ws6.Range("B19").Value = "0.31 (-0.15, 0.8)"

$ws6.Range("A20").Value = "Middle school"
$ws6.Range("B20").Value = "0.27 (-0.19, 0.72)"

$ws6.Range("A21").Value = "Change in county COVID-19 case rate"
$ws6.Range("B21").Value = "-0.08 (-0.45, 0.3)"

# ---------------------------------------------------------------------------
# Sheet tbl7 (sheet7): rows reshuffled/rewritten + new rows 5-14
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("tbl7")

$ws7.Range("B2").Value = "1.07 (0.33, 1.76)"

$ws7.Range("A3").Value = "Physical distancing"
$ws7.Range("B3").Value = "-0.2 (-0.73, 0.31)"

$ws7.Range("A4").Value = "Staying home when sick"
$ws7.Range("B4").Value = "-0.09 (-0.67, 0.39)"

$ws7.Range("A5").Value = "HVAC systems"
$ws7.Range("B5").Value = "-0.21 (-0.82, 0.36)"

$ws7.Range("A6").Value = "Percent two or more races"
$ws7.Range("B6").Value = "0.25 (0.04, 0.46)"

$ws7.Range("A7").Value = "Percent Asian"
$ws7.Range("B7").Value = "0.01 (-0.2, 0.21)"

$ws7.Range("A8").Value = "Percent White"
$ws7.Range("B8").Value = "-0.19 (-0.52, 0.15)"

$ws7.Range("A9").Value = "Percent free and reduced lunch"
$ws7.Range("B9").Value = "-0.28 (-0.6, 0.03)"

$ws7.Range("A10").Value = "SVI Overall Rank"
$ws7.Range("B10").Value = "0.03 (-0.18, 0.26)"

$ws7.Range("A11").Value = "Percent Black or African American"
$ws7.Range("B11").Value = "-0.08 (-0.33, 0.21)"

$ws7.Range("A12").Value = "High school"
$ws7.Range("B12").Value = "0.3 (-0.17, 0.8)"

$ws7.Range("A13").Value = "Middle school"
$ws7.Range("B13").Value = "0.27 (-0.21, 0.73)"

$ws7.Range("A14").Value = "Change in county COVID-19 case rate"
$ws7.Range("B14").Value = "-0.09 (-0.45, 0.22)"

# ---------------------------------------------------------------------------
# New sheet tbl8 (sheet8): appended after tbl7
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "tbl8"

$ws8.Range("A1").Value = "Strategy"
$ws8.Range("B1").Value = "Coefficient (95% interval)"
$ws8.Range("A1:B1").Font.Bold = $true
$ws8.Range("A1:B1").HorizontalAlignment = -4108

$ws8.Range("A2").Value = "Intercept"
$ws8.Range("B2").Value = "1.02 (0.34, 1.67)"

$ws8.Range("A3").Value = "One strategy"
$ws8.Range("B3").Value = "0.15 (-0.41, 0.77)"

$ws8.Range("A4").Value = "Two strategies"
$ws8.Range("B4").Value = "-0.23 (-0.79, 0.29)"

$ws8.Range("A5").Value = "Three strategies"
$ws8.Range("B5").Value = "-0.53 (-1.07, -0.05)"

$ws8.Range("A6").Value = "Percent two or more races"
$ws8.Range("B6").Value = "0.24 (0.03, 0.46)"

$ws8.Range("A7").Value = "Percent Asian"
$ws8.Range("B7").Value = "0.01 (-0.22, 0.22)"

$ws8.Range("A8").Value = "Percent White"
$ws8.Range("B8").Value = "-0.19 (-0.55, 0.19)"

$ws8.Range("A9").Value = "Percent free and reduced lunch"
$ws8.Range("B9").Value = "-0.28 (-0.58, 0.05)"

$ws8.Range("A10").Value = "SVI Overall Rank"
$ws8.Range("B10").Value = "0.03 (-0.2, 0.27)"

$ws8.Range("A11").Value = "Percent Black or African American"
$ws8.Range("B11").Value = "-0.07 (-0.35, 0.16)"

$ws8.Range("A12").Value = "High school"
$ws8.Range("B12").Value = "0.28 (-0.18, 0.74)"

$ws8.Range("A13").Value = "Middle school"
$ws8.Range("B13").Value = "0.24 (-0.18, 0.71)"

$ws8.Range("A14").Value = "Change in county COVID-19 case rate"
$ws8.Range("B14").Value = "-0.09 (-0.41, 0.3)"
